# Commit: "updated date on slide"
#   - Slide 1: update the date textbox from "Monday" to "Monday 26, 2026".
#   - Slide 25: two small wording fixes in the timer-explanation body text
#     ("At "+"200 MHz:" -> "At 200 MHz:" and "thus to ensure" -> "thus, to ensure").
$p = $ppt.ActivePresentation

# --- Slide 1: date textbox ---
$s1 = $p.Slides.Item(1)
$dateShape = $s1.Shapes.Item("Google Shape;83;p26")
$dateShape.TextFrame.TextRange.Text = "Monday 26, 2026"

# --- Slide 25: body placeholder text ---
$s25 = $p.Slides.Item(25)
$bodyShape = $s25.Shapes.Item("Google Shape;118;p29")
$bodyTr = $bodyShape.TextFrame.TextRange

# "At " + "200 MHz:" -> single run "At 200 MHz:"
$found1 = $bodyTr.Find("At 200 MHz:")
$found1.Text = "At 200 MHz:"

# "thus to ensure" -> "thus, to ensure"
$found2 = $bodyTr.Find("s, thus to ensure that your timer runs every 10 ")
$found2.Text = "s, thus, to ensure that your timer runs every 10 "
